# Apply "updated heat transfer models" edit to user_inputs.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sensitivity_variables")

# Update the mach-limit baseline value for row 27 (B27)
$ws.Range("B27").Value = 0.3

# Replace the "Y" placeholder values in C41:C43 and C46:C48 with
# numeric heat-transfer model coefficients
$ws.Range("C41").Value = 3
$ws.Range("C42").Value = 2.5
$ws.Range("C43").Value = 2
$ws.Range("C46").Value = 2
$ws.Range("C47").Value = 2.5
$ws.Range("C48").Value = 3

# Update the active/selected cell and scroll position on the sheet view
$ws.Activate()
$ws.Range("C49").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
